# Feat: Implementar y depurar plataforma de castings funcional
#
# Adds a new casting entry ("Alba Galocha", Costa Rica) as row 4 of the
# "Mujeres" sheet, matching the formatting already used by the existing
# data rows, and widens columns A:M on that sheet to a uniform width.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mujeres")

# Carry over the formatting from the last existing data row (row 3) onto
# the new row 4 before writing values into it.
$ws.Range("A3:M3").Copy()
$ws.Range("A4:M4").PasteSpecial(-4122)

$ws.Range("A4").Value = "Alba"
$ws.Range("B4").Value = "Galocha"
$ws.Range("C4").Value = "Costa Rica"
$ws.Range("D4").Value = 178.0
$ws.Range("E4").Value = 88.0
$ws.Range("F4").Value = 64.0
$ws.Range("G4").Value = 92.0
$ws.Range("H4").Value = 39.0
$ws.Range("I4").Value = "Rubio"
$ws.Range("J4").Value = "Azules"
$ws.Range("K4").Value = "valenrojas"
$ws.Range("L4").Value = "valenrojastk"
$ws.Range("M4").Value = "alba-galocha"

# Widen columns A:M uniformly (closest attainable value to the target
# 17.63 character-width given Excel's pixel-grid rounding of ColumnWidth).
$ws.Range("A1:M4").Columns.ColumnWidth = 16.76
